$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) column cells are plain numeric-looking strings in the source data but
# must remain stored as TEXT (matching the original inlineStr cells), so force the
# format to Text before assigning, then drop back to the Normal style so no stray
# number-format style lingers on the cell.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.774.62'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.681.64'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.85%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.19'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.96%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.547'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.681.07'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.85%  '
$ws.Range('E10').Value = '  +2.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.159'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.364'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.24'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.99'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.158.30'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.20%  '
$ws.Range('E16').Value = '  -2.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.798.79'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.685.11'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.78'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.92'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '365.59'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.34%  '
$ws.Range('E22').Value = '  -3.21%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.84'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.11%  '
$ws.Range('E24').Value = '  -3.85%  '
$ws.Range('E25').Value = '  +0.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '70.99'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.94%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.23'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.834.38'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('E29').Value = '  -2.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '555.49'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.03'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.39%  '
$ws.Range('E33').Value = '  -3.60%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.94'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.16%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.131'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.23%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('E37').Value = '  -4.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.54'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.61%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '154.61'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.94%  '
$ws.Range('E40').Value = '  -2.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.32'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.57%  '
$ws.Range('E42').Value = '  -3.97%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.94'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.39%  '
$ws.Range('E44').Value = '  -6.02%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.46'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0₆0302'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.592'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '154.14'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.77%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.89'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.68%  '
$ws.Range('E51').Value = '  -3.45%  '
